$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# remain text (matching the source data's inlineStr/text convention) instead
# of being auto-converted to a number by Excel's input parser. We do this by
# temporarily marking the cell as Text (@) before assigning the value, then
# resetting the cell style back to Normal so no stray number format sticks to
# the cell (it still leaves the value stored as text).
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D18', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D44', 'D45', 'D46', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume change, and for a few rows
# the coin name/link swapped positions with a neighboring row).
$ws.Range('D2').Value = '69.806.51'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '3.526.72'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '607.98'
$ws.Range('E5').Value = '  +3.11%  '
$ws.Range('D6').Value = '184.89'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '3.523.28'
$ws.Range('E7').Value = '  -1.20%  '
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '0.211'
$ws.Range('E10').Value = '  +5.43%  '
$ws.Range('D11').Value = '0.639'
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').Value = '53.39'
$ws.Range('E12').Value = '  -2.55%  '
$ws.Range('D13').Value = '0.0000306'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').Value = '4.088.29'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '69.836.04'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.550.49'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '580.47'
$ws.Range('E18').Value = '  +3.82%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '12.57'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '18.81'
$ws.Range('E20').Value = '  -3.51%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '0.987'
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('D23').Value = '17.31'
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('D24').Value = '4.68'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '4.82'
$ws.Range('E25').Value = '  -2.38%  '
$ws.Range('D26').Value = '93.64'
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '2.95'
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '10.94'
$ws.Range('E28').Value = '  -4.92%  '
$ws.Range('D29').Value = '9.33'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').Value = '31.94'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('D31').Value = '6.98'
$ws.Range('D32').Value = '12.13'
$ws.Range('E32').Value = '  -3.38%  '
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('D34').Value = '63.02'
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('D35').Value = '3.30'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').Value = '3.58'
$ws.Range('E36').Value = '  +16.90%  '
$ws.Range('D37').Value = '531.96'
$ws.Range('E37').Value = '  -4.44%  '
$ws.Range('D38').Value = '0.401'
$ws.Range('E38').Value = '  -3.44%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').Value = '36.96'
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('D41').Value = '0.0₃0776'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').Value = '3.515.77'
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').Value = '0.135'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').Value = '0.0452'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').Value = '3.44'
$ws.Range('E46').Value = '  -4.40%  '
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('D49').Value = '9.03'
$ws.Range('E49').Value = '  -2.95%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').Value = '1.42'
$ws.Range('E51').Value = '  -1.48%  '

# Reset the number format on the forced-text cells back to Normal so no
# leftover style attribute is left on the cell itself.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
